# Updated cryptos list values (Price and Volume(1h) columns), per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.248.19"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "'1.786.36"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'334.60"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.3799"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3417"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'48.48"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "'1.197"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "'0.07493"
$ws.Range("E11").Value = "  -3.13%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'21.89"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "'6.464"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "'1.786.43"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "'0.06650"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "'83.88"
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'6.627"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").Value = "'17.33"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "'27.248.31"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  -5.73%  "
$ws.Range("D25").Value = "'2.410"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "'1.494"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'2.542"
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").Value = "'21.32"
$ws.Range("E28").Value = "  -4.16%  "
$ws.Range("D29").Value = "'153.15"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'1.989.22"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "'134.27"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'4.022"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").Value = "'6.079"
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "'13.31"
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").Value = "'1.653"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").Value = "'0.6955"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").Value = "'5.448"
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("D39").Value = "'0.2210"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "'8.826"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "'0.06331"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").Value = "'0.02337"
$ws.Range("E42").Value = "  -3.22%  "
$ws.Range("D43").Value = "'1.236"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "'14.38"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("D45").Value = "'0.6516"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'3.842"
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("D48").Value = "'2.151"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "'129.37"
$ws.Range("D50").Value = "'0.07130"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").Value = "'79.02"
$ws.Range("E51").Value = "  -1.75%  "
